$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "0189e477-237d-494f-ab96-2e21c6669357"
$ws.Range("B3").Value = "João"
$ws.Range("C3").Value = "dautonico32@gmail.com"

$ws.Range("A4").Value = "3d7aecf9-7827-425e-a099-4bd3f943153c"
$ws.Range("B4").Value = "Funcionário Base"
$ws.Range("C4").Value = "base@empresa.com"
